$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 2024
$ws.Range("B2").Value = "DEC"
$ws.Range("C2").Value = "31/12-01/12"
$ws.Range("D2").Value = "-"
